# Adds four new benchmark rows (v17-2200, v18-500600, v18-510800,
# v18-519300) to the "Sheet1" benchmark table and turns on 3-colour-scale
# conditional formatting for the Win-Rate / Score (average) / Time (avg)
# columns, matching the "Various benchmarks of V18" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New rows -------------------------------------------------------
# Write the text cells first, and in this particular order, so the
# resulting shared-string table is built up in the same sequence the
# original author typed them in (A20, then C21, then A21, A22, A23).
$ws.Range("A20").Value = "v17-2200"
$ws.Range("C21").Value = "Added power pellets to state"
$ws.Range("A21").Value = "v18-500600"
$ws.Range("A22").Value = "v18-510800"
$ws.Range("A23").Value = "v18-519300"

# Row 20 - v17-2200
$ws.Range("B20").Value = 5325
$ws.Range("D20").Value = 176
$ws.Range("E20").Value = 0.78
$ws.Range("F20").Value = 2.09
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 3931.59
$ws.Range("J20").Value = 1630
$ws.Range("K20").Value = 7600
$ws.Range("L20").Value = 239.8
$ws.Range("M20").Value = 135
$ws.Range("N20").Value = 244
$ws.Range("O20").Value = 4.48
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 9
$ws.Range("R20").Value = 113.44
$ws.Range("S20").Value = 56.3
$ws.Range("T20").Value = 194.2

# Row 21 - v18-500600
$ws.Range("B21").Value = 11395
$ws.Range("D21").Value = 394
$ws.Range("E21").Value = 0.44
$ws.Range("F21").Value = 1.15
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 3159.39
$ws.Range("J21").Value = 1880
$ws.Range("K21").Value = 5450
$ws.Range("L21").Value = 235.48
$ws.Range("M21").Value = 174
$ws.Range("N21").Value = 244
$ws.Range("O21").Value = 2.47
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 7
$ws.Range("R21").Value = 118.17
$ws.Range("S21").Value = 59
$ws.Range("T21").Value = 253

# Row 22 - v18-510800
$ws.Range("B22").Value = 11415
$ws.Range("D22").Value = 385
$ws.Range("E22").Value = 0.62
$ws.Range("F22").Value = 1.76
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 5
$ws.Range("I22").Value = 4225.66
$ws.Range("J22").Value = 2340
$ws.Range("K22").Value = 9800
$ws.Range("L22").Value = 239.92
$ws.Range("M22").Value = 164
$ws.Range("N22").Value = 244
$ws.Range("O22").Value = 4.98
$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 12
$ws.Range("R22").Value = 123.52
$ws.Range("S22").Value = 59.1
$ws.Range("T22").Value = 245.8

# Row 23 - v18-519300
$ws.Range("B23").Value = 11422
$ws.Range("D23").Value = 188
$ws.Range("E23").Value = 0.85
$ws.Range("F23").Value = 2.39
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 3980.32
$ws.Range("J23").Value = 2800
$ws.Range("K23").Value = 7400
$ws.Range("L23").Value = 242.65
$ws.Range("M23").Value = 218
$ws.Range("N23").Value = 244
$ws.Range("O23").Value = 4.56
$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 9
$ws.Range("R23").Value = 101.71
$ws.Range("S23").Value = 64.5
$ws.Range("T23").Value = 207.8

# --- Conditional formatting ------------------------------------------
# Three colour-scale rules over whole columns I (Score average),
# E (Win-Rate) and R (Time avg). They are added in this order so the
# resulting priorities come out as I=3, E=2, R=1 (last-added rule gets
# top priority), matching the target workbook.
$fcI = $ws.Range("I1:I1048576").FormatConditions.AddColorScale(3)
$fcE = $ws.Range("E1:E1048576").FormatConditions.AddColorScale(3)
$fcR = $ws.Range("R1:R1048576").FormatConditions.AddColorScale(3)
$fcI.Priority = 3
$fcE.Priority = 2
$fcR.Priority = 1

# --- Selection ---------------------------------------------------------
# After entering the new rows, the user's active cell moved to A24.
$ws.Range("A24").Select()

"Inserted rows 20-23 and added colour-scale conditional formatting"
